$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.062
$ws.Range("E2").Value = 0.021
$ws.Range("F2").Value = 0.018
$ws.Range("H2").Value = 0.039
$ws.Range("I2").Value = 0.033
$ws.Range("J2").Value = 0.034

# Row 3
$ws.Range("C3").Value = 0.083
$ws.Range("D3").Value = 0.071
$ws.Range("E3").Value = 0.065
$ws.Range("F3").Value = 0.051
$ws.Range("G3").Value = 0.103
$ws.Range("H3").Value = 0.111
$ws.Range("I3").Value = 0.068
$ws.Range("J3").Value = 0.081

# Row 4
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("C4").Value = 0.021
$ws.Range("D4").Value = 0.041
$ws.Range("E4").Value = 0.044
$ws.Range("F4").Value = 0.033
$ws.Range("G4").Value = 0.067
$ws.Range("H4").Value = 0.072
$ws.Range("I4").Value = 0.035

# Row 5
$ws.Range("B5").Value = 0.579
$ws.Range("C5").Value = 0.104
$ws.Range("D5").Value = 0.078
$ws.Range("E5").Value = 0.066
$ws.Range("F5").Value = 0.048
$ws.Range("G5").Value = 0.093
$ws.Range("H5").Value = 0.096
$ws.Range("I5").Value = 0.074
$ws.Range("J5").Value = 0.081

# Row 6
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.579
$ws.Range("C6").Value = 0.042
$ws.Range("D6").Value = 0.048
$ws.Range("E6").Value = 0.045
$ws.Range("F6").Value = 0.03
$ws.Range("G6").Value = 0.057
$ws.Range("H6").Value = 0.057
$ws.Range("I6").Value = 0.041
$ws.Range("J6").Value = 0.047

# Row 7
$ws.Range("B7").Value = 0.579
$ws.Range("C7").Value = 0.021
$ws.Range("D7").Value = 0.007
$ws.Range("E7").Value = 0.001
$ws.Range("F7").Value = -0.003
$ws.Range("G7").Value = -0.01
$ws.Range("H7").Value = -0.015
$ws.Range("I7").Value = 0.007
$ws.Range("J7").Value = 0
